$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "55.688.18"
Set-TextValue "E2" "  -2.14%  "
Set-TextValue "D3" "2.966.17"
Set-TextValue "E3" "  -4.37%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  -0.12%  "
Set-TextValue "D5" "492.33"
Set-TextValue "E5" "  -4.42%  "
Set-TextValue "D6" "132.73"
Set-TextValue "E6" "  +1.38%  "
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  -0.26%  "
Set-TextValue "D8" "2.958.23"
Set-TextValue "E8" "  -4.57%  "
Set-TextValue "D10" "7.12"
Set-TextValue "E10" "  +0.95%  "
Set-TextValue "E11" "  -5.11%  "
Set-TextValue "D12" "0.350"
Set-TextValue "E12" "  -6.86%  "
Set-TextValue "E13" "  +1.16%  "
Set-TextValue "D14" "3.476.61"
Set-TextValue "E14" "  -5.31%  "
Set-TextValue "D15" "24.93"
Set-TextValue "E15" "  -0.69%  "
Set-TextValue "D16" "55.594.01"
Set-TextValue "E16" "  -2.59%  "
Set-TextValue "D17" "2.960.31"
Set-TextValue "E17" "  -5.04%  "
Set-TextValue "E18" "  -4.00%  "
Set-TextValue "D19" "5.71"
Set-TextValue "E19" "  +0.69%  "
Set-TextValue "D20" "12.20"
Set-TextValue "E20" "  -4.05%  "
Set-TextValue "D21" "7.51"
Set-TextValue "E21" "  -3.88%  "
Set-TextValue "D22" "321.30"
Set-TextValue "E22" "  -5.74%  "
Set-TextValue "E23" "  -0.25%  "
Set-TextValue "D24" "0.465"
Set-TextValue "E24" "  -6.40%  "
Set-TextValue "D25" "60.39"
Set-TextValue "E25" "  -10.91%  "
Set-TextValue "D26" "1.00"
Set-TextValue "E26" "  +1.27%  "
Set-TextValue "E27" "  +1.80%  "
Set-TextValue "E28" "  +0.02%  "
Set-TextValue "D29" "0.0₃0849"
Set-TextValue "E29" "  -7.19%  "
Set-TextValue "D30" "6.60"
Set-TextValue "E30" "  +0.51%  "
Set-TextValue "D31" "6.58"
Set-TextValue "E31" "  -3.09%  "
Set-TextValue "D32" "1.18"
Set-TextValue "E32" "  +0.69%  "
Set-TextValue "D33" "1.71"
Set-TextValue "E33" "  -6.38%  "
Set-TextValue "D34" "19.54"
Set-TextValue "E34" "  -8.21%  "
Set-TextValue "D35" "149.20"
Set-TextValue "E35" "  -4.31%  "
Set-TextValue "D36" "4.41"
Set-TextValue "E36" "  -6.63%  "
Set-TextValue "E37" "  -4.19%  "
Set-TextValue "D38" "5.70"
Set-TextValue "E38" "  -6.16%  "
Set-TextValue "D39" "23.35"
Set-TextValue "E39" "  -8.54%  "
Set-TextValue "D40" "0.0652"
Set-TextValue "E40" "  -3.63%  "
Set-TextValue "D41" "2.996.52"
Set-TextValue "E41" "  -5.38%  "
Set-TextValue "D42" "1.00"
Set-TextValue "E42" "  -0.18%  "
Set-TextValue "D43" "36.34"
Set-TextValue "E43" "  -9.36%  "
Set-TextValue "D44" "1.00"
Set-TextValue "E44" "  -2.37%  "
Set-TextValue "D45" "0.631"
Set-TextValue "E45" "  -6.62%  "
Set-TextValue "E46" "  -3.38%  "
Set-TextValue "D47" "3.54"
Set-TextValue "E47" "  -7.42%  "
Set-TextValue "D48" "2.134.10"
Set-TextValue "E48" "  -3.87%  "
Set-TextValue "D49" "0.0236"
Set-TextValue "E49" "  +2.54%  "
Set-TextValue "D50" "19.39"
Set-TextValue "E50" "  -1.24%  "
Set-TextValue "D51" "5.60"
Set-TextValue "E51" "  -7.65%  "
